# Update countries & provincias Spain
# - Refresh the "last updated" timestamp.
# - Re-sort a handful of countries whose updated case counts changed their
#   rank in the (descending, by total cases) table, and write the refreshed
#   case numbers for the affected countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp (row 1) -----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 12:22"

# --- Rows 34/35: Mexico <-> Rumania swap places, Rumania gets new data -----
$ws.Range("A34").Value = "Rumania"
$ws.Range("B34").Value = 9710
$ws.Range("C34").Value = 468
$ws.Range("D34").Value = 2406
$ws.Range("E34").Value = 6796
$ws.Range("F34").Value = 288
$ws.Range("G34").Value = 10
$ws.Range("H34").Value = 508

$ws.Range("A35").Value = "Mexico"
$ws.Range("B35").Value = 9501
$ws.Range("C35").Value = 729
$ws.Range("D35").Value = 2627
$ws.Range("E35").Value = 6017
$ws.Range("F35").Value = 378
$ws.Range("G35").Value = 145
$ws.Range("H35").Value = 857

# --- Row 44: Australia, numbers refreshed only ------------------------------
$ws.Range("B44").Value = 6649
$ws.Range("C44").Value = 4
$ws.Range("E44").Value = 1655

# --- Row 56: Marruecos, numbers refreshed only ------------------------------
$ws.Range("B56").Value = 3377
$ws.Range("C56").Value = 168
$ws.Range("D56").Value = 398
$ws.Range("E56").Value = 2830
$ws.Range("G56").Value = 4
$ws.Range("H56").Value = 149

# --- Rows 76/77: Eslovenia <-> Bosnia y Herzegovina swap --------------------
$ws.Range("A76").Value = "Bosnia y Herzegovina"
$ws.Range("B76").Value = 1368
$ws.Range("C76").Value = 26
$ws.Range("D76").Value = 460
$ws.Range("E76").Value = 855
$ws.Range("F76").Value = 4
$ws.Range("H76").Value = 53

$ws.Range("A77").Value = "Eslovenia"
$ws.Range("B77").Value = 1353
$ws.Range("C77").Value = 9
$ws.Range("D77").Value = 205
$ws.Range("E77").Value = 1069
$ws.Range("F77").Value = 24
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = 79

# --- Rows 98/99/100: Albania moves up above Kirguistan/Bolivia -------------
$ws.Range("A98").Value = "Albania"
$ws.Range("B98").Value = 634
$ws.Range("C98").Value = 25
$ws.Range("D98").Value = 356
$ws.Range("E98").Value = 251
$ws.Range("F98").Value = 4
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 27

$ws.Range("A99").Value = "Kirguistan"
$ws.Range("B99").Value = 612
$ws.Range("C99").Value = 22
$ws.Range("D99").Value = 254
$ws.Range("E99").Value = 351
$ws.Range("F99").Value = 5
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 7

$ws.Range("A100").Value = "Bolivia"
$ws.Range("C100").Value = 11
$ws.Range("D100").Value = 44
$ws.Range("E100").Value = 528
$ws.Range("F100").Value = 3
$ws.Range("G100").Value = 3
$ws.Range("H100").Value = 37

# --- Row 135: Brunei, numbers refreshed only --------------------------------
$ws.Range("D135").Value = 117
$ws.Range("E135").Value = 20

# --- Rows 140/141: Trinidad y Tobago <-> Etiopia swap -----------------------
$ws.Range("A140").Value = "Etiopia"
$ws.Range("B140").Value = 116
$ws.Range("C140").Value = 2
$ws.Range("D140").Value = 21
$ws.Range("E140").Value = 92
$ws.Range("H140").Value = 3

$ws.Range("A141").Value = "Trinidad yTobago"
$ws.Range("B141").Value = 115
$ws.Range("D141").Value = 28
$ws.Range("E141").Value = 79
$ws.Range("H141").Value = 8

# --- Row 152: Zambia, numbers refreshed only --------------------------------
$ws.Range("B152").Value = 74
$ws.Range("C152").Value = 4
$ws.Range("E152").Value = 36
